# Minor updates to lectures 2 and 4:
# remove the trailing blank slide (slide 13) from the deck.
$p = $ppt.ActivePresentation
$p.Slides.Item($p.Slides.Count).Delete()
